# "blok 10 userstories van de website af"
# Mark a block of userstories on the "Basis" sheet as finished:
# - A14 -> "Kom er niet uit" (stuck / couldn't figure it out)
# - A15, A16, A17, A18 -> "Afgerond" (done)
# Also move the visible scroll position / selection down a bit,
# from A14 to A17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basis")

$ws.Range("A14").Value = "Kom er niet uit"
$ws.Range("A15").Value = "Afgerond"
$ws.Range("A16").Value = "Afgerond"
$ws.Range("A17").Value = "Afgerond"
$ws.Range("A18").Value = "Afgerond"

$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A17").Select() | Out-Null
